$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "max" column (column C) entirely; this shifts
# "prediction" (D) -> C and "rejection-f" (E) -> D.
$ws.Range("C1:C3").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

# Update the numeric prediction-difference values in column B.
$ws.Range("B2").Value = 0.0871431380098926
$ws.Range("B3").Value = -0.04862773867910164

# Mark row 3's rejection-f value as rejected.
$ws.Range("D3").Value = "s__CAG-631 sp000433015(reject)"
